# Apply scheduled-runner updates to Gilgamesh_Profits workbook
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 250012400
$ws.Range("I86").Value = 400018400
$ws.Range("J86").Value = 2399.3333
$ws.Range("K86").Value = 400018400
$ws.Range("L86").Value = 2399.3333
$ws.Range("M86").Value = -400017277
$ws.Range("N86").Value = -4645.3333
$ws.Range("H89").Value = 250012400
$ws.Range("I89").Value = 400018400
$ws.Range("J89").Value = 2399.3333
$ws.Range("K89").Value = 2000092000
$ws.Range("L89").Value = 11996.6665
$ws.Range("M89").Value = -2000086384
$ws.Range("N89").Value = -23228.6665
$ws.Range("H113").Value = 47621108
$ws.Range("J113").Value = 2400
$ws.Range("L113").Value = 2400
$ws.Range("N113").Value = -8908
$ws.Range("H137").Value = 3539
$ws.Range("I137").Value = 3223.8333
$ws.Range("J137").Value = 4011.75
$ws.Range("K137").Value = 9671.499899999999
$ws.Range("L137").Value = 12035.25
$ws.Range("M137").Value = -7121.499899999999
$ws.Range("N137").Value = -17135.25
$ws.Range("H138").Value = 398080.2
$ws.Range("J138").Value = 446283.16
$ws.Range("L138").Value = 1338849.48
$ws.Range("N138").Value = -1349129.48

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 37931.12
$ws.Range("I32").Value = 11327.607
$ws.Range("K32").Value = 11327.607
$ws.Range("M32").Value = -11040.607
$ws.Range("H61").Value = 6422.5884
$ws.Range("I61").Value = 2463.2222
$ws.Range("K61").Value = 2463.2222
$ws.Range("M61").Value = -2251.2222
$ws.Range("H74").Value = 152644.45
$ws.Range("I74").Value = 207699.38
$ws.Range("J74").Value = 3996.2
$ws.Range("K74").Value = 207699.38
$ws.Range("L74").Value = 3996.2
$ws.Range("M74").Value = -206825.38
$ws.Range("N74").Value = -5744.2
$ws.Range("H77").Value = 152644.45
$ws.Range("I77").Value = 207699.38
$ws.Range("J77").Value = 3996.2
$ws.Range("K77").Value = 1038496.9
$ws.Range("L77").Value = 19981
$ws.Range("M77").Value = -1034128.9
$ws.Range("N77").Value = -28717
$ws.Range("H132").Value = 2028
$ws.Range("I132").Value = 1642.7097
$ws.Range("K132").Value = 4928.1291
$ws.Range("M132").Value = -2398.1291
$ws.Range("H134").Value = 89726.5
$ws.Range("J134").Value = 89726.5
$ws.Range("L134").Value = 89726.5
$ws.Range("N134").Value = -99866.5
$ws.Range("H135").Value = 107498
$ws.Range("J135").Value = 107498
$ws.Range("L135").Value = 107498
$ws.Range("N135").Value = -117638
$ws.Range("H136").Value = 6422.5884
$ws.Range("I136").Value = 2463.2222
$ws.Range("K136").Value = 7389.6666
$ws.Range("M136").Value = -4839.6666
$ws.Range("H138").Value = 102834
$ws.Range("J138").Value = 102834
$ws.Range("L138").Value = 102834
$ws.Range("N138").Value = -113114

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3492.5
$ws.Range("I86").Value = 2740.75
$ws.Range("K86").Value = 2740.75
$ws.Range("M86").Value = -1617.75
$ws.Range("H89").Value = 3492.5
$ws.Range("I89").Value = 2740.75
$ws.Range("K89").Value = 13703.75
$ws.Range("M89").Value = -8087.75
$ws.Range("H105").Value = 5910985.5
$ws.Range("I105").Value = 304732.66
$ws.Range("K105").Value = 304732.66
$ws.Range("M105").Value = -302985.66
$ws.Range("H134").Value = 2802.975
$ws.Range("I134").Value = 2623.25
$ws.Range("K134").Value = 7869.75
$ws.Range("M134").Value = -5334.75
$ws.Range("H135").Value = 115661.43
$ws.Range("J135").Value = 115661.43
$ws.Range("L135").Value = 115661.43
$ws.Range("N135").Value = -125801.43
$ws.Range("H141").Value = 78000
$ws.Range("J141").Value = 78000
$ws.Range("L141").Value = 78000
$ws.Range("N141").Value = -88360

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3715.3794
$ws.Range("J31").Value = 4093
$ws.Range("L31").Value = 4093
$ws.Range("N31").Value = -4683
$ws.Range("H34").Value = 3715.3794
$ws.Range("J34").Value = 4093
$ws.Range("L34").Value = 4093
$ws.Range("N34").Value = -4497
$ws.Range("H58").Value = 4008.762
$ws.Range("I58").Value = 3727.5
$ws.Range("J58").Value = 4571.2856
$ws.Range("K58").Value = 3727.5
$ws.Range("L58").Value = 4571.2856
$ws.Range("M58").Value = -3524.5
$ws.Range("N58").Value = -4977.2856
$ws.Range("H62").Value = 14296899
$ws.Range("I62").Value = 25004076
$ws.Range("J62").Value = 20665
$ws.Range("K62").Value = 25004076
$ws.Range("L62").Value = 20665
$ws.Range("M62").Value = -25003452
$ws.Range("N62").Value = -21913
$ws.Range("H65").Value = 14296899
$ws.Range("I65").Value = 25004076
$ws.Range("J65").Value = 20665
$ws.Range("K65").Value = 125020380
$ws.Range("L65").Value = 103325
$ws.Range("M65").Value = -125017260
$ws.Range("N65").Value = -109565
$ws.Range("H99").Value = 5360
$ws.Range("I99").Value = 5446.778
$ws.Range("K99").Value = 5446.778
$ws.Range("M99").Value = -3948.778
$ws.Range("H126").Value = 5360
$ws.Range("I126").Value = 5446.778
$ws.Range("K126").Value = 16340.334
$ws.Range("M126").Value = -13870.334
$ws.Range("H134").Value = 2918.8108
$ws.Range("I134").Value = 2078.0312
$ws.Range("K134").Value = 6234.0936
$ws.Range("M134").Value = -3699.0936
$ws.Range("H136").Value = 4008.762
$ws.Range("I136").Value = 3727.5
$ws.Range("J136").Value = 4571.2856
$ws.Range("K136").Value = 11182.5
$ws.Range("L136").Value = 13713.8568
$ws.Range("M136").Value = -8632.5
$ws.Range("N136").Value = -18813.8568

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 158.3
$ws.Range("J23").Value = 266.4
$ws.Range("L23").Value = 799.1999999999999
$ws.Range("N23").Value = -1269.2
$ws.Range("H33").Value = 288.22726
$ws.Range("I33").Value = 140.8
$ws.Range("K33").Value = 844.8000000000001
$ws.Range("M33").Value = -561.8000000000001
$ws.Range("H40").Value = 355.4
$ws.Range("I40").Value = 265.08334
$ws.Range("J40").Value = 716.6667
$ws.Range("K40").Value = 1060.33336
$ws.Range("L40").Value = 2866.6668
$ws.Range("M40").Value = -991.3333600000001
$ws.Range("N40").Value = -3004.6668
$ws.Range("H107").Value = 1652.6
$ws.Range("I107").Value = 1800
$ws.Range("J107").Value = 1636.2222
$ws.Range("K107").Value = 5400
$ws.Range("L107").Value = 4908.6666
$ws.Range("M107").Value = -3480
$ws.Range("N107").Value = -8748.6666
$ws.Range("H122").Value = 2176
$ws.Range("J122").Value = 2476.2856
$ws.Range("L122").Value = 22286.5704
$ws.Range("N122").Value = -27186.5704
$ws.Range("H139").Value = 3114.25
$ws.Range("I139").Value = 2392.3845
$ws.Range("J139").Value = 3739.8667
$ws.Range("K139").Value = 7177.1535
$ws.Range("L139").Value = 11219.6001
$ws.Range("M139").Value = -2037.1535
$ws.Range("N139").Value = -21499.6001

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 111113336
$ws.Range("I80").Value = 333334660
$ws.Range("K80").Value = 333334660
$ws.Range("M80").Value = -333333662
$ws.Range("H83").Value = 111113336
$ws.Range("I83").Value = 333334660
$ws.Range("K83").Value = 1666673300
$ws.Range("M83").Value = -1666668308
$ws.Range("H132").Value = 2564.7334
$ws.Range("I132").Value = 1684.6522
$ws.Range("K132").Value = 5053.9566
$ws.Range("M132").Value = -2523.9566
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6326.8438
$ws.Range("I122").Value = 6249.35
$ws.Range("K122").Value = 18748.05
$ws.Range("M122").Value = -16298.05
$ws.Range("H136").Value = 6062.25
$ws.Range("I136").Value = 3499.6667
$ws.Range("J136").Value = 13750
$ws.Range("K136").Value = 10499.0001
$ws.Range("L136").Value = 41250
$ws.Range("M136").Value = -7949.000100000001
$ws.Range("N136").Value = -46350

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6287.0977
$ws.Range("I132").Value = 5936.3145
$ws.Range("K132").Value = 17808.9435
$ws.Range("M132").Value = -15278.9435
$ws.Range("H138").Value = 115997.2
$ws.Range("J138").Value = 115997.2
$ws.Range("L138").Value = 115997.2
$ws.Range("N138").Value = -126277.2
